$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 3 columns: A=query, B=dbExcel, C=WebExcel (headers)
# and on row 2: A=big Neo4j query, B=Neo4jData filename, C=WebData filename.
# We need to insert a new column between A and B that holds the new
# "StatQuery" header/value (stat-bar count query), shifting the old B/C to C/D.

# Insert a new column before column B - this shifts the existing B -> C and C -> D
$ws.Columns("B:B").Insert()

# Match the new column's width to column A's width (same as the other "query" column)
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# New header for the inserted column
$ws.Range("B1").Value = "StatQuery"

# New stat-bar query text (wraps, like the big query in A2) for the inserted column
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Irish Wolfhound']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match the row styling so B2 wraps text the same way A2 does
$ws.Range("B2").WrapText = $true

# Move the active selection onto the new query cell
$ws.Range("B2").Select()
